$wb = $excel.ActiveWorkbook

$wsPmax   = $wb.Worksheets.Item("Pmax_Pgen.prn")
$wsMayor  = $wb.Worksheets.Item("Mayor_maxima.prn")
$wsMenor  = $wb.Worksheets.Item("Menor_optima.prn")

# --- Header relabeling (shared across the three report sheets) ---
# "RESERVA[%]" -> "RESERVA_DATO[%]"
$wsPmax.Range("G1").Value  = "RESERVA_DATO[%]"
$wsMayor.Range("G1").Value = "RESERVA_DATO[%]"
$wsMenor.Range("G1").Value = "RESERVA_DATO[%]"

# Pmax_Pgen.prn: I1 "DATO" -> "RES_OPT[%]", drop the extra J1 ("RES_OPT[[%]") column
$wsPmax.Range("I1").Value = "RES_OPT[%]"
$wsPmax.Range("J1").ClearContents()

# Mayor_maxima.prn / Menor_optima.prn: I1 "DATO" -> "RESOPT[%]" (absorbs old J1 label), drop J1
$wsMayor.Range("I1").Value = "RESOPT[%]"
$wsMayor.Range("J1").ClearContents()

$wsMenor.Range("I1").Value = "RESOPT[%]"
$wsMenor.Range("J1").ClearContents()

# --- New data rows (informes al 60 %) for Pmax_Pgen.prn and Mayor_maxima.prn ---
$rows = @(
    @(101, "NUC-A       21.600", 1, 944.9999570846558, 749.9989624023438, 195.000994682312,  26.0001685945936,  5,  3),
    @(102, "NUC-B       21.600", 1, 944.9999570846558, 749.9989624023438, 195.000994682312,  26.0001685945936,  5,  0),
    @(206, "URBGEN      18.000", 1, 899.9999761581421, 799.9994506835938, 100.0005254745483, 12.50007426743839, 5,  1),
    @(211, "HYDRO_G     20.000", 1, 800.4000333607197, 579.9990844726562, 220.4009488880635, 38.00022358456915, 10, 0)
)

foreach ($ws in @($wsPmax, $wsMayor)) {
    $r = 2
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
        $ws.Cells.Item($r, 7).Value = $row[6]
        $ws.Cells.Item($r, 8).Value = $row[7]
        $ws.Cells.Item($r, 9).Value = $row[8]
        $r++
    }
}
